$d = $word.ActiveDocument

# 1. Update the date
$d.Content.Find.Execute("2025-02-05", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-13", 2)

# 2. Fix missing space: "важливіше сфокусуватися"
$d.Content.Find.Execute("важливішесфокусуватися", $true, $false, $false, $false, $false, $true, 1, $false, "важливіше сфокусуватися", 2)

# 3. Remove "Винахідник мережі Інтернет " prefix
$d.Content.Find.Execute("Винахідник мережі Інтернет Тім Бернерс-Лі", $true, $false, $false, $false, $false, $true, 1, $false, "Тім Бернерс-Лі", 2)

# 4. Swap timestamps between row 4 and row 5 in the verbatim block
$d.Content.Find.Execute(" 4 application/rdf+xml 2018-11-05T10:21:52.515944", $true, $false, $false, $false, $false, $true, 1, $false, " 4 application/rdf+xml 2018-11-05T10:37:34.390529", 2)
$d.Content.Find.Execute(" 5 application/rdf+xml 2018-11-05T10:37:34.390529", $true, $false, $false, $false, $false, $true, 1, $false, " 5 application/rdf+xml 2018-11-05T10:21:52.515944", 2)
